$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume snapshot (scheduled GitHub Actions update).
# Column D ("Price") cells whose new text is a plain decimal number are written with a
# leading apostrophe so Excel keeps them as literal text (preserving trailing zeros and
# thousand-grouped dotted values like "3.542.24") instead of auto-converting to Number.
$ws.Range("D2").Value = "65.379.12"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "3.542.24"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'599.60"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").Value = "'135.79"
$ws.Range("E6").Value = "  -1.63%  "
$ws.Range("D7").Value = "3.543.35"
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").Value = "'0.497"
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("D11").Value = "'7.17"
$ws.Range("E11").Value = "  +3.35%  "
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("D13").Value = "4.145.62"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").Value = "'27.63"
$ws.Range("E14").Value = "  +0.87%  "
$ws.Range("D15").Value = "'0.0000183"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").Value = "3.548.10"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").Value = "65.407.39"
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("D19").Value = "'9.84"
$ws.Range("E19").Value = "  -2.61%  "
$ws.Range("D20").Value = "'14.55"
$ws.Range("E20").Value = "  +2.30%  "
$ws.Range("D21").Value = "'5.77"
$ws.Range("E21").Value = "  -1.86%  "
$ws.Range("D22").Value = "'394.20"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").Value = "'0.582"
$ws.Range("E23").Value = "  +1.27%  "
$ws.Range("D24").Value = "3.687.26"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "'74.71"
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "'0.0000116"
$ws.Range("E27").Value = "  +2.04%  "
$ws.Range("D28").Value = "'7.93"
$ws.Range("E28").Value = "  +1.21%  "
$ws.Range("D29").Value = "'1.63"
$ws.Range("E29").Value = "  +15.54%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'8.53"
$ws.Range("E31").Value = "  +2.74%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'2.31"
$ws.Range("E32").Value = "  +1.63%  "
$ws.Range("D33").Value = "3.550.65"
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").Value = "'24.29"
$ws.Range("E34").Value = "  +1.61%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  +1.60%  "
$ws.Range("D37").Value = "'5.33"
$ws.Range("E37").Value = "  +5.83%  "
$ws.Range("D38").Value = "'1.59"
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("D39").Value = "'169.15"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").Value = "'6.91"
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("D41").Value = "'0.0830"
$ws.Range("E41").Value = "  +2.98%  "
$ws.Range("D42").Value = "'0.830"
$ws.Range("E42").Value = "  +0.59%  "
$ws.Range("D43").Value = "'26.22"
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("D44").Value = "'1.26"
$ws.Range("E44").Value = "  +4.33%  "
$ws.Range("D45").Value = "'42.89"
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("D46").Value = "'0.999"
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "'4.48"
$ws.Range("E47").Value = "  +0.65%  "
$ws.Range("D48").Value = "'1.68"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("D49").Value = "'6.93"
$ws.Range("E49").Value = "  +1.15%  "
$ws.Range("D50").Value = "2.391.43"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").Value = "'0.899"
$ws.Range("E51").Value = "  +5.63%  "
